# #5: property aircraft done
# Fix the "property_category" column on the 建物 (building) and 汽車 (car)
# sheets: both were incorrectly populated with the shared string "land"
# (copy/paste leftover from the 土地 sheet). Correct them to "building"
# and "car" respectively.

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet - property_category column I, data rows 2-5
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I5").Value = "building"

# 汽車 (Car) sheet - property_category column H, data row 2
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
